$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 299.25
$ws.Range("I18").Value = 299.25
$ws.Range("K18").Value = 299.25
$ws.Range("M18").Value = -15.25

$ws.Range("H28").Value = 2227.8
$ws.Range("I28").Value = 3047.5
$ws.Range("J28").Value = 1681.3334
$ws.Range("K28").Value = 3047.5
$ws.Range("L28").Value = 1681.3334
$ws.Range("M28").Value = -2562.5
$ws.Range("N28").Value = -2651.3334

$ws.Range("H43").Value = 27424.273
$ws.Range("J43").Value = 21441.8
$ws.Range("L43").Value = 21441.8
$ws.Range("N43").Value = -21579.8

$ws.Range("H54").Value = 5341.1665
$ws.Range("I54").Value = 5341.1665
$ws.Range("K54").Value = 5341.1665
$ws.Range("M54").Value = -4855.1665

$ws.Range("H62").Value = 6045.1465
$ws.Range("I62").Value = 3666.1428
$ws.Range("K62").Value = 3666.1428
$ws.Range("M62").Value = -3042.1428

$ws.Range("H65").Value = 6045.1465
$ws.Range("I65").Value = 3666.1428
$ws.Range("K65").Value = 18330.714
$ws.Range("M65").Value = -15210.714

$ws.Range("H74").Value = 92677.86
$ws.Range("I74").Value = 157570.14
$ws.Range("K74").Value = 157570.14
$ws.Range("M74").Value = -156634.14

$ws.Range("H77").Value = 92677.86
$ws.Range("I77").Value = 157570.14
$ws.Range("K77").Value = 787850.7000000001
$ws.Range("M77").Value = -783170.7000000001

$ws.Range("H116").Value = 12560.625
$ws.Range("I116").Value = 15121.25
$ws.Range("K116").Value = 15121.25
$ws.Range("M116").Value = -11679.25

$ws.Range("H141").Value = 1693
$ws.Range("I141").Value = 1693
$ws.Range("K141").Value = 5079
$ws.Range("M141").Value = 101

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5548.9
$ws.Range("I2").Value = 4638.4287
$ws.Range("K2").Value = 4638.4287
$ws.Range("M2").Value = -4525.4287

$ws.Range("H60").Value = 45749.5
$ws.Range("I60").Value = 45749.5
$ws.Range("K60").Value = 45749.5
$ws.Range("M60").Value = -45016.5

$ws.Range("H116").Value = 5548.9
$ws.Range("I116").Value = 4638.4287
$ws.Range("K116").Value = 4638.4287
$ws.Range("M116").Value = -2344.4287

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5548.9
$ws.Range("I3").Value = 4638.4287
$ws.Range("K3").Value = 4638.4287
$ws.Range("M3").Value = -4524.4287

$ws.Range("H29").Value = 5995.2856
$ws.Range("J29").Value = 3129.5
$ws.Range("L29").Value = 3129.5
$ws.Range("N29").Value = -3707.5

$ws.Range("H36").Value = 15042.286
$ws.Range("I36").Value = 882.6667
$ws.Range("K36").Value = 882.6667
$ws.Range("M36").Value = -348.6667

$ws.Range("H86").Value = 2846
$ws.Range("I86").Value = 2799.5
$ws.Range("J86").Value = 3001
$ws.Range("K86").Value = 2799.5
$ws.Range("L86").Value = 3001
$ws.Range("M86").Value = -1676.5
$ws.Range("N86").Value = -5247

$ws.Range("H89").Value = 2846
$ws.Range("I89").Value = 2799.5
$ws.Range("J89").Value = 3001
$ws.Range("K89").Value = 13997.5
$ws.Range("L89").Value = 15005
$ws.Range("M89").Value = -8381.5
$ws.Range("N89").Value = -26237

$ws.Range("H94").Value = 2096.6956
$ws.Range("I94").Value = 1749.238
$ws.Range("J94").Value = 5745
$ws.Range("K94").Value = 1749.238
$ws.Range("L94").Value = 5745
$ws.Range("M94").Value = -1298.238
$ws.Range("N94").Value = -6647

$ws.Range("H99").Value = 4510.5
$ws.Range("I99").Value = 1378.1428
$ws.Range("J99").Value = 7642.857
$ws.Range("K99").Value = 1378.1428
$ws.Range("L99").Value = 7642.857
$ws.Range("M99").Value = 119.8571999999999
$ws.Range("N99").Value = -10638.857

$ws.Range("H134").Value = 1990.15
$ws.Range("I134").Value = 1716.8334
$ws.Range("K134").Value = 5150.5002
$ws.Range("M134").Value = -2615.5002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2170
$ws.Range("I105").Value = 1750
$ws.Range("K105").Value = 1750
$ws.Range("M105").Value = -3

$ws.Range("H134").Value = 12342.03
$ws.Range("I134").Value = 8938.821
$ws.Range("K134").Value = 26816.463
$ws.Range("M134").Value = -24281.463

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 56868.945
$ws.Range("J131").Value = 2313
$ws.Range("L131").Value = 6939
$ws.Range("N131").Value = -17019

$ws.Range("H134").Value = 5650.913
$ws.Range("I134").Value = 4103.737
$ws.Range("K134").Value = 12311.211
$ws.Range("M134").Value = -7241.210999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -330
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4075.125
$ws.Range("I61").Value = 3999.7693
$ws.Range("K61").Value = 3999.7693
$ws.Range("M61").Value = -3797.7693

$ws.Range("H68").Value = 6138.125
$ws.Range("I68").Value = 4100.8335
$ws.Range("J68").Value = 12250
$ws.Range("K68").Value = 4100.8335
$ws.Range("L68").Value = 12250
$ws.Range("M68").Value = -3351.8335
$ws.Range("N68").Value = -13748

$ws.Range("H71").Value = 6138.125
$ws.Range("I71").Value = 4100.8335
$ws.Range("J71").Value = 12250
$ws.Range("K71").Value = 20504.1675
$ws.Range("L71").Value = 61250
$ws.Range("M71").Value = -16760.1675
$ws.Range("N71").Value = -68738

$ws.Range("H93").Value = 4624.8
$ws.Range("I93").Value = 4305.3335
$ws.Range("J93").Value = 7500
$ws.Range("K93").Value = 4305.3335
$ws.Range("L93").Value = 7500
$ws.Range("M93").Value = -3057.3335
$ws.Range("N93").Value = -9996

$ws.Range("H101").Value = 10866.5
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 10866.5
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 10866.5
$ws.Range("N101").Value = -17356.5
$ws.Range("M101").ClearContents()

$ws.Range("H113").Value = 4075.125
$ws.Range("I113").Value = 3999.7693
$ws.Range("K113").Value = 3999.7693
$ws.Range("M113").Value = -1829.7693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2330
$ws.Range("I96").Value = 2216.6667
$ws.Range("J96").Value = 2500
$ws.Range("K96").Value = 2216.6667
$ws.Range("L96").Value = 2500
$ws.Range("M96").Value = -843.6667000000002
$ws.Range("N96").Value = -5246

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
